# Auto-generated Excel COM edit script
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 53.46154
$ws.Range("I4").Value = 40
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 40
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = 74
$ws.Range("H5").Value = 34.2
$ws.Range("I5").Value = 34.2
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 34.2
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 80.8
$ws.Range("H6").Value = 414
$ws.Range("I6").Value = 394
$ws.Range("J6").Value = 494
$ws.Range("K6").Value = 1182
$ws.Range("L6").Value = 1482
$ws.Range("M6").Value = -1070
$ws.Range("H15").Value = 524.9773
$ws.Range("I15").Value = 524.9773
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1574.9319
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1405.9319
$ws.Range("H18").Value = 437.8889
$ws.Range("I18").Value = 437.8889
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 437.8889
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -153.8889
$ws.Range("H38").Value = 4347.08
$ws.Range("I38").Value = 2476.3076
$ws.Range("J38").Value = 6373.75
$ws.Range("K38").Value = 7428.9228
$ws.Range("L38").Value = 19121.25
$ws.Range("M38").Value = -7056.9228
$ws.Range("N38").Value = -19865.25
$ws.Range("H41").Value = 147.66667
$ws.Range("I41").Value = 85.875
$ws.Range("J41").Value = 271.25
$ws.Range("K41").Value = 85.875
$ws.Range("L41").Value = 271.25
$ws.Range("M41").Value = 354.125
$ws.Range("H43").Value = 2077.3333
$ws.Range("I43").Value = 1904.5
$ws.Range("J43").Value = 2215.6
$ws.Range("K43").Value = 1904.5
$ws.Range("L43").Value = 2215.6
$ws.Range("M43").Value = -1835.5
$ws.Range("N43").Value = -2353.6
$ws.Range("H51").Value = 4970.619
$ws.Range("I51").Value = 3400
$ws.Range("J51").Value = 5049.15
$ws.Range("K51").Value = 3400
$ws.Range("L51").Value = 5049.15
$ws.Range("M51").Value = -2916
$ws.Range("N51").Value = -6017.15
$ws.Range("H58").Value = 1216.1428
$ws.Range("I58").Value = 380.75
$ws.Range("J58").Value = 2330
$ws.Range("K58").Value = 1142.25
$ws.Range("L58").Value = 6990
$ws.Range("M58").Value = -992.25
$ws.Range("N58").Value = -7290
$ws.Range("H87").Value = 84474.5
$ws.Range("I87").Value = 78999
$ws.Range("J87").Value = 89950
$ws.Range("K87").Value = 78999
$ws.Range("L87").Value = 89950
$ws.Range("M87").Value = -77751
$ws.Range("N87").Value = -92446
$ws.Range("H90").Value = 84474.5
$ws.Range("I90").Value = 78999
$ws.Range("J90").Value = 89950
$ws.Range("K90").Value = 236997
$ws.Range("L90").Value = 269850
$ws.Range("M90").Value = -230757
$ws.Range("N90").Value = -282330
$ws.Range("H113").Value = 9717.556
$ws.Range("I113").Value = 8126.25
$ws.Range("J113").Value = 10990.6
$ws.Range("K113").Value = 8126.25
$ws.Range("L113").Value = 10990.6
$ws.Range("M113").Value = -4872.25
$ws.Range("H115").Value = 286.625
$ws.Range("I115").Value = 286.625
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 859.875
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 707.125
$ws.Range("H129").Value = 14680.125
$ws.Range("I129").Value = 1309.4286
$ws.Range("J129").Value = 25079.555
$ws.Range("K129").Value = 3928.2858
$ws.Range("L129").Value = 75238.66500000001
$ws.Range("M129").Value = 1071.7142
$ws.Range("H138").Value = 3021.923
$ws.Range("I138").Value = 2444
$ws.Range("J138").Value = 3070.0833
$ws.Range("K138").Value = 7332
$ws.Range("L138").Value = 9210.249899999999
$ws.Range("M138").Value = -2192
$ws.Range("N138").Value = -19490.2499

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3717.276
$ws.Range("I32").Value = 3335.2593
$ws.Range("J32").Value = 8874.5
$ws.Range("K32").Value = 3335.2593
$ws.Range("L32").Value = 8874.5
$ws.Range("M32").Value = -3048.2593
$ws.Range("H45").Value = 3502.75
$ws.Range("I45").Value = 2004
$ws.Range("J45").Value = 7999
$ws.Range("K45").Value = 2004
$ws.Range("L45").Value = 7999
$ws.Range("M45").Value = -1627
$ws.Range("N45").Value = -8753
$ws.Range("H102").Value = 289428.84
$ws.Range("I102").Value = 336167
$ws.Range("J102").Value = 9000
$ws.Range("K102").Value = 336167
$ws.Range("L102").Value = 9000
$ws.Range("M102").Value = -334545
$ws.Range("H110").Value = 117214.336
$ws.Range("I110").Value = 117214.336
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 117214.336
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -115169.336
$ws.Range("H122").Value = 7753856
$ws.Range("I122").Value = 1743.6552
$ws.Range("J122").Value = 23811804
$ws.Range("K122").Value = 5230.9656
$ws.Range("L122").Value = 71435412
$ws.Range("M122").Value = -2780.9656
$ws.Range("H132").Value = 24433014
$ws.Range("I132").Value = 9839.764999999999
$ws.Range("J132").Value = 143059860
$ws.Range("K132").Value = 29519.295
$ws.Range("L132").Value = 429179580
$ws.Range("M132").Value = -26989.295
$ws.Range("N132").Value = -429184640
$ws.Range("H135").Value = 104852.125
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 104852.125
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 104852.125
$ws.Range("N135").Value = -114992.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 525.2857
$ws.Range("I22").Value = 381.875
$ws.Range("J22").Value = 716.5
$ws.Range("K22").Value = 381.875
$ws.Range("L22").Value = 716.5
$ws.Range("M22").Value = -208.875
$ws.Range("N22").Value = -1062.5
$ws.Range("H135").Value = 82333
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 82333
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 82333
$ws.Range("N135").Value = -92473

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("H108").Value = 29536.8
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 29536.8
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 29536.8
$ws.Range("N108").Value = -37216.8
$ws.Range("H132").Value = 67356.58
$ws.Range("I132").Value = 69555.13
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 208665.39
$ws.Range("L132").Value = 4200
$ws.Range("M132").Value = -206135.39
$ws.Range("N132").Value = -9260
$ws.Range("M50").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 583.3333
$ws.Range("I16").Value = 425
$ws.Range("J16").Value = 662.5
$ws.Range("K16").Value = 1275
$ws.Range("L16").Value = 1987.5
$ws.Range("M16").Value = -1102
$ws.Range("N16").Value = -2333.5
$ws.Range("H38").Value = 54.916668
$ws.Range("I38").Value = 14
$ws.Range("J38").Value = 177.66667
$ws.Range("K38").Value = 42
$ws.Range("L38").Value = 533.00001
$ws.Range("M38").Value = 305
$ws.Range("H40").Value = 154
$ws.Range("I40").Value = 58.857143
$ws.Range("J40").Value = 265
$ws.Range("K40").Value = 235.428572
$ws.Range("L40").Value = 1060
$ws.Range("M40").Value = -166.428572
$ws.Range("N40").Value = -1198
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("H129").Value = 11495317
$ws.Range("I129").Value = 22222786
$ws.Range("J129").Value = 1600.7142
$ws.Range("K129").Value = 66668358
$ws.Range("L129").Value = 4802.142599999999
$ws.Range("M129").Value = -66663358
$ws.Range("H131").Value = 23810836
$ws.Range("I131").Value = 41667610
$ws.Range("J131").Value = 1810.5555
$ws.Range("K131").Value = 125002830
$ws.Range("L131").Value = 5431.666499999999
$ws.Range("M131").Value = -124997790
$ws.Range("N131").Value = -15511.6665
$ws.Range("N68").ClearContents()
$ws.Range("N71").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 26317122
$ws.Range("I122").Value = 1270.75
$ws.Range("J122").Value = 166668340
$ws.Range("K122").Value = 3812.25
$ws.Range("L122").Value = 500005020
$ws.Range("M122").Value = -1362.25
$ws.Range("N122").Value = -500009920

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1200
$ws.Range("I55").Value = 834.5
$ws.Range("J55").Value = 1565.5
$ws.Range("K55").Value = 834.5
$ws.Range("L55").Value = 1565.5
$ws.Range("M55").Value = -661.5
$ws.Range("N55").Value = -1911.5
$ws.Range("H132").Value = 4971.108
$ws.Range("I132").Value = 2613.2693
$ws.Range("J132").Value = 10544.182
$ws.Range("K132").Value = 7839.8079
$ws.Range("L132").Value = 31632.546
$ws.Range("M132").Value = -5309.8079
$ws.Range("N132").Value = -36692.546
$ws.Range("H136").Value = 2962.3684
$ws.Range("I136").Value = 1952.6666
$ws.Range("J136").Value = 6748.75
$ws.Range("K136").Value = 5857.9998
$ws.Range("L136").Value = 20246.25
$ws.Range("M136").Value = -3307.9998
$ws.Range("N136").Value = -25346.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1429583.9
$ws.Range("I100").Value = 1429583.9
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2859167.8
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2858626.8
$ws.Range("H136").Value = 11226.583
$ws.Range("I136").Value = 4681.091
$ws.Range("J136").Value = 16765.076
$ws.Range("K136").Value = 14043.273
$ws.Range("L136").Value = 50295.228
$ws.Range("M136").Value = -11493.273
$ws.Range("N136").Value = -55395.228
$ws.Range("H141").Value = 116330
$ws.Range("I141").Value = 109000
$ws.Range("J141").Value = 119995
$ws.Range("K141").Value = 109000
$ws.Range("L141").Value = 119995
$ws.Range("M141").Value = -103820
$ws.Range("N141").Value = -130355
